# "form for staff added" — adds a BATCH column (L) to the student marks sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text fix: REGNO -> REG.NO. -----------------------------------
$ws.Range("A1").Value = "REG.NO."

# --- New BATCH column (L) --------------------------------------------------
# 1) Copy K1's current "last column" look (right border style s=4) onto the
#    new last column L1 *before* touching K1's own formatting.
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)

# 2) K1 is no longer the last column, so it reverts to the regular header
#    formatting (same as the rest of the header row, e.g. J1).
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)

# 3) Give the new data cells L2:L40 the same formatting as the other mark
#    columns (style s=3, matching column C for example).
$ws.Range("C2:C40").Copy()
$ws.Range("L2:L40").PasteSpecial(-4122)

# --- Values ------------------------------------------------------------
$ws.Range("L1").Value = "BATCH"

$ws.Range("L2").Value  = "2015 - 2018 "
$ws.Range("L3").Value  = "2015 - 2018 "
$ws.Range("L4").Value  = "2015 - 2018"
$ws.Range("L5").Value  = "2015 - 2018 "
$ws.Range("L6").Value  = "2015 - 2018"
$ws.Range("L7").Value  = "2015 - 2018"
$ws.Range("L8").Value  = "2015 - 2018"
$ws.Range("L9").Value  = "2015 - 2018"
$ws.Range("L10").Value = "2015 - 2018"
$ws.Range("L11").Value = "2015 - 2018"
$ws.Range("L12").Value = "2015 - 2018"
$ws.Range("L13").Value = "2015 - 2018"
$ws.Range("L14").Value = "2015 - 2018"
$ws.Range("L15").Value = "2015 - 2018"
$ws.Range("L16").Value = "2015 - 2018"
$ws.Range("L17").Value = "2015 - 2018"
$ws.Range("L18").Value = "2015 - 2018"
$ws.Range("L19").Value = "2015 - 2018"
$ws.Range("L20").Value = "2015 - 2018"
$ws.Range("L21").Value = "2015 - 2018 "
$ws.Range("L22").Value = "2015 - 2018"
$ws.Range("L23").Value = "2015 - 2018"
$ws.Range("L24").Value = "2015 - 2018"
$ws.Range("L25").Value = "2015 - 2018"
$ws.Range("L26").Value = "2015 - 2018"
$ws.Range("L27").Value = "2015 - 2018 "
$ws.Range("L28").Value = "2015 - 2018 "
$ws.Range("L29").Value = "2015 - 2018"
$ws.Range("L30").Value = "2015 - 2018 "
$ws.Range("L31").Value = "2015 - 2018"
$ws.Range("L32").Value = "2015 - 2018"
$ws.Range("L33").Value = "2015 - 2018"
$ws.Range("L34").Value = "2015 - 2018"
$ws.Range("L35").Value = "2015 - 2018"
$ws.Range("L36").Value = "2015 - 2018"
$ws.Range("L37").Value = "2015 - 2018"
$ws.Range("L38").Value = "2015 - 2018"
$ws.Range("L39").Value = "2015 - 2018"
$ws.Range("L40").Value = "2015 - 2018"

# --- Conditional formatting: extend "= u" red-bold rule to the new column --
$fc = $ws.Range("C3").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("C3:L40"))

Write-Host "BATCH column added"
